$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 0.001
$ws.Range("C2").Value = 0.8300114457558273
$ws.Range("D2").Value = 0.7915470841965732
$ws.Range("E2").Value = 0.7782804609958788
$ws.Range("F2").Value = 0.8284575478508657
$ws.Range("G2").Value = 0.7688617347966722
$ws.Range("H2").Value = 0.8240723542800704
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 0.8295297742443334
$ws.Range("D3").Value = 0.792894407632566
$ws.Range("E3").Value = 0.7783343005057132
$ws.Range("F3").Value = 0.8284319002284807
$ws.Range("G3").Value = 0.7685955222220381
$ws.Range("H3").Value = 0.8238437536823762
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 0.2
$ws.Range("C4").Value = 0.8273510348240765
$ws.Range("D4").Value = 0.791679895033229
$ws.Range("E4").Value = 0.7761615612038371
$ws.Range("F4").Value = 0.8271270088148999
$ws.Range("G4").Value = 0.7671546042467117
$ws.Range("H4").Value = 0.8225830749308044
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 0.3
$ws.Range("C5").Value = 0.8289023745289795
$ws.Range("D5").Value = 0.7911050863825863
$ws.Range("E5").Value = 0.7768582104647008
$ws.Range("F5").Value = 0.8253892955251706
$ws.Range("G5").Value = 0.7645872318682709
$ws.Range("H5").Value = 0.8207604875019869
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0.4
$ws.Range("C6").Value = 0.8274348343094843
$ws.Range("D6").Value = 0.7873699205600815
$ws.Range("E6").Value = 0.77371019819639
$ws.Range("F6").Value = 0.8215022894017578
$ws.Range("G6").Value = 0.7597771217327933
$ws.Range("H6").Value = 0.8167724325280509
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.8290481388278373
$ws.Range("D7").Value = 0.78637269753646
$ws.Range("E7").Value = 0.7733760753431626
$ws.Range("F7").Value = 0.8164934976011895
$ws.Range("G7").Value = 0.7540640662914666
$ws.Range("H7").Value = 0.8116982198262392
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 0.6
$ws.Range("C8").Value = 0.8252139663441143
$ws.Range("D8").Value = 0.7852555317118037
$ws.Range("E8").Value = 0.7709409303718634
$ws.Range("F8").Value = 0.8133768768790195
$ws.Range("G8").Value = 0.750097821778277
$ws.Range("H8").Value = 0.8084242215279149
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 0.7
$ws.Range("C9").Value = 0.824502889615484
$ws.Range("D9").Value = 0.7811185446376095
$ws.Range("E9").Value = 0.7677309600574654
$ws.Range("F9").Value = 0.8089445195074247
$ws.Range("G9").Value = 0.7446130268843639
$ws.Range("H9").Value = 0.8039031926447576
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 0.8
$ws.Range("C10").Value = 0.8242550800931426
$ws.Range("D10").Value = 0.7794412229886776
$ws.Range("E10").Value = 0.7662984979372558
$ws.Range("F10").Value = 0.8038453457556874
$ws.Range("G10").Value = 0.738022375388897
$ws.Range("H10").Value = 0.7984995282374352
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = 0.9
$ws.Range("C11").Value = 0.8281017219503934
$ws.Range("D11").Value = 0.7780193962011592
$ws.Range("E11").Value = 0.7672883547216045
$ws.Range("F11").Value = 0.803398763279219
$ws.Range("G11").Value = 0.7364718684501106
$ws.Range("H11").Value = 0.797887671440335
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.8257809839657323
$ws.Range("D12").Value = 0.7745863463029883
$ws.Range("E12").Value = 0.7637978155767071
$ws.Range("F12").Value = 0.7974591623991957
$ws.Range("G12").Value = 0.7298635632700952
$ws.Range("H12").Value = 0.7917898360985201
